$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final outcome measures added in column D (rows 2-15), in row order
$values = @(
    "Not worse",        # row 2
    "A little worse",   # row 3
    "A little worse",   # row 4
    "Not worse",         # row 5
    "A little worse",   # row 6
    "A little worse",   # row 7
    "Not worse",         # row 8
    "Somewhat worse",   # row 9
    "Somewhat worse",   # row 10
    "Not worse",         # row 11
    "Not worse",         # row 12
    "A little worse",   # row 13
    "A little worse",   # row 14
    "Not worse"          # row 15
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Widen column D slightly to fit the new content
$ws.Columns.Item(4).ColumnWidth = 18.2825

# Update the active selection to reflect the new last-used cell
$ws.Range("D16").Select()
